# Auto-generated Excel COM-interop script
# Adds 2024-08-30 violent-crime data: updates column K (2024 totals)
# across Citywide Totals, By Neighborhood, and each neighborhood sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 11).Value = 5256
$ws.Cells.Item(3, 11).Value = 5426
$ws.Cells.Item(4, 11).Value = 1125
$ws.Cells.Item(6, 11).Value = 6037
$ws.Cells.Item(7, 11).Value = 18231

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(5, 11).Value = 42
$ws.Cells.Item(6, 11).Value = 133
$ws.Cells.Item(7, 11).Value = 538
$ws.Cells.Item(8, 11).Value = 1230
$ws.Cells.Item(11, 11).Value = 349
$ws.Cells.Item(18, 11).Value = 123
$ws.Cells.Item(19, 11).Value = 536
$ws.Cells.Item(21, 11).Value = 57
$ws.Cells.Item(22, 11).Value = 48
$ws.Cells.Item(23, 11).Value = 191
$ws.Cells.Item(27, 11).Value = 173
$ws.Cells.Item(29, 11).Value = 976
$ws.Cells.Item(31, 11).Value = 199
$ws.Cells.Item(33, 11).Value = 779
$ws.Cells.Item(36, 11).Value = 241
$ws.Cells.Item(37, 11).Value = 614
$ws.Cells.Item(40, 11).Value = 42
$ws.Cells.Item(41, 11).Value = 127
$ws.Cells.Item(42, 11).Value = 674
$ws.Cells.Item(43, 11).Value = 161
$ws.Cells.Item(45, 11).Value = 21
$ws.Cells.Item(47, 11).Value = 123
$ws.Cells.Item(48, 11).Value = 228
$ws.Cells.Item(51, 11).Value = 229
$ws.Cells.Item(53, 11).Value = 235
$ws.Cells.Item(54, 11).Value = 358
$ws.Cells.Item(55, 11).Value = 204
$ws.Cells.Item(57, 11).Value = 67
$ws.Cells.Item(60, 11).Value = 113
$ws.Cells.Item(63, 11).Value = 50
$ws.Cells.Item(64, 11).Value = 116
$ws.Cells.Item(65, 11).Value = 415
$ws.Cells.Item(67, 11).Value = 692
$ws.Cells.Item(73, 11).Value = 156
$ws.Cells.Item(76, 11).Value = 252
$ws.Cells.Item(78, 11).Value = 209
$ws.Cells.Item(79, 11).Value = 449
$ws.Cells.Item(83, 11).Value = 405
$ws.Cells.Item(85, 11).Value = 856
$ws.Cells.Item(86, 11).Value = 123
$ws.Cells.Item(90, 11).Value = 164
$ws.Cells.Item(91, 11).Value = 198
$ws.Cells.Item(94, 11).Value = 240
$ws.Cells.Item(95, 11).Value = 312
$ws.Cells.Item(96, 11).Value = 198
$ws.Cells.Item(101, 11).Value = 18231

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(6, 11).Value = 86
$ws.Cells.Item(7, 11).Value = 198

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 11).Value = 183
$ws.Cells.Item(7, 11).Value = 538

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(3, 11).Value = 91
$ws.Cells.Item(7, 11).Value = 349

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 11).Value = 286
$ws.Cells.Item(7, 11).Value = 856

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(3, 11).Value = 60
$ws.Cells.Item(7, 11).Value = 235

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 11).Value = 337
$ws.Cells.Item(6, 11).Value = 418
$ws.Cells.Item(7, 11).Value = 1230

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(3, 11).Value = 148
$ws.Cells.Item(7, 11).Value = 405

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 11).Value = 214
$ws.Cells.Item(3, 11).Value = 288
$ws.Cells.Item(6, 11).Value = 225
$ws.Cells.Item(7, 11).Value = 779

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(2, 11).Value = 104
$ws.Cells.Item(7, 11).Value = 312

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 11).Value = 176
$ws.Cells.Item(3, 11).Value = 204
$ws.Cells.Item(4, 11).Value = 29
$ws.Cells.Item(6, 11).Value = 178
$ws.Cells.Item(7, 11).Value = 614

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(2, 11).Value = 131
$ws.Cells.Item(7, 11).Value = 415

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(3, 11).Value = 48
$ws.Cells.Item(7, 11).Value = 199

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(4, 11).Value = 39
$ws.Cells.Item(7, 11).Value = 692

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(2, 11).Value = 57
$ws.Cells.Item(6, 11).Value = 191
$ws.Cells.Item(7, 11).Value = 358

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 11).Value = 279
$ws.Cells.Item(3, 11).Value = 352
$ws.Cells.Item(6, 11).Value = 270
$ws.Cells.Item(7, 11).Value = 976

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(2, 11).Value = 30
$ws.Cells.Item(3, 11).Value = 53
$ws.Cells.Item(7, 11).Value = 228

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(6, 11).Value = 171
$ws.Cells.Item(7, 11).Value = 536

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(6, 11).Value = 136
$ws.Cells.Item(7, 11).Value = 252

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Cells.Item(3, 11).Value = 37
$ws.Cells.Item(7, 11).Value = 133

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Cells.Item(6, 11).Value = 49
$ws.Cells.Item(7, 11).Value = 127

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 11).Value = 178
$ws.Cells.Item(3, 11).Value = 209
$ws.Cells.Item(6, 11).Value = 254
$ws.Cells.Item(7, 11).Value = 674

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(2, 11).Value = 62
$ws.Cells.Item(7, 11).Value = 209

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(3, 11).Value = 58
$ws.Cells.Item(7, 11).Value = 204

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(6, 11).Value = 52
$ws.Cells.Item(7, 11).Value = 191

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(3, 11).Value = 94
$ws.Cells.Item(7, 11).Value = 198

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Cells.Item(6, 11).Value = 34
$ws.Cells.Item(7, 11).Value = 57

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(2, 11).Value = 150
$ws.Cells.Item(7, 11).Value = 449

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Cells.Item(2, 11).Value = 22
$ws.Cells.Item(6, 11).Value = 44
$ws.Cells.Item(7, 11).Value = 116

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(3, 11).Value = 38
$ws.Cells.Item(7, 11).Value = 123

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(2, 11).Value = 94
$ws.Cells.Item(6, 11).Value = 55
$ws.Cells.Item(7, 11).Value = 241

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(3, 11).Value = 46
$ws.Cells.Item(7, 11).Value = 240

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(2, 11).Value = 35
$ws.Cells.Item(3, 11).Value = 36
$ws.Cells.Item(7, 11).Value = 123

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(2, 11).Value = 50
$ws.Cells.Item(7, 11).Value = 156

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Cells.Item(2, 11).Value = 10
$ws.Cells.Item(7, 11).Value = 42

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(2, 11).Value = 47
$ws.Cells.Item(7, 11).Value = 173

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Cells.Item(6, 11).Value = 30
$ws.Cells.Item(7, 11).Value = 123

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(2, 11).Value = 58
$ws.Cells.Item(7, 11).Value = 164

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(2, 11).Value = 66
$ws.Cells.Item(7, 11).Value = 229

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Cells.Item(3, 11).Value = 13
$ws.Cells.Item(7, 11).Value = 67

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Cells.Item(4, 11).Value = 9
$ws.Cells.Item(7, 11).Value = 113

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(4, 11).Value = 22
$ws.Cells.Item(7, 11).Value = 161

$ws = $wb.Worksheets.Item('Clearing')
$ws.Cells.Item(3, 11).Value = 16
$ws.Cells.Item(7, 11).Value = 48

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Cells.Item(6, 11).Value = 10
$ws.Cells.Item(7, 11).Value = 21

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(7, 11).Value = 42

Write-Host "Updated 154 cells"